# Update the variant calling lecture and lab:
# Remove the two "Annovar" slides (gene information / allele frequency)
# that are no longer part of the lecture. After each delete, the slides
# below shift up one position, so we repeatedly remove slide index 2.

$p = $ppt.ActivePresentation

$p.Slides.Item(2).Delete()   # was "Annovar and gene information"
$p.Slides.Item(2).Delete()   # was "Annovar and allele frequency"
